# Update Keras Hyperband Optimizer
# Move the "Max"/"Min" values (columns D/E) into "Mean"/"Standard Deviation"
# (columns B/C) for rows 2-4 on the DOE sheet, and clear out the old D/E
# cells entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DOE")

# Row 2 - POX/C
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 10
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 - C/A
$ws.Range("B3").Value = 0.1
$ws.Range("C3").Value = 0.01
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

# Row 4 - POX/M
$ws.Range("B4").Value = 0.001
$ws.Range("C4").Value = 0.0001
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
